# "finished size up room analysis"
# Re-run of the timing numbers (and a few experimental counts) on the
# "increasing room size" sheet, plus the resulting tab/selection switch.

$wb = $excel.ActiveWorkbook
$wsRoom = $wb.Worksheets.Item("room constraint")
$ws = $wb.Worksheets.Item("increasing room size")

# --- Updated "Time (s)" column (G2:G16) -----------------------------------
$ws.Range("G2").Value = 0.0130233764648437
$ws.Range("G3").Value = 0.010282039642333899
$ws.Range("G4").Value = 0.010424852371215799
$ws.Range("G5").Value = 0.0131230354309082
$ws.Range("G6").Value = 0.0098819732666015608
$ws.Range("G7").Value = 0.011027812957763601
$ws.Range("G8").Value = 0.012421846389770499
$ws.Range("G9").Value = 0.011552095413207999
$ws.Range("G10").Value = 0.0126750469207763
$ws.Range("G11").Value = 0.0107879638671875
$ws.Range("G12").Value = 0.0160486698150634
$ws.Range("G13").Value = 0.012868881225585899
$ws.Range("G14").Value = 0.012700080871582
$ws.Range("G15").Value = 0.016459941864013599
$ws.Range("G16").Value = 0.014331102371215799

# --- Updated "Experimental" counts (I column) for the rows that moved -----
$ws.Range("I3").Value = 3173
$ws.Range("I5").Value = 3161
$ws.Range("I16").Value = 4169

# --- Page setup: force portrait orientation on "increasing room size" -----
$ws.PageSetup.Orientation = 1

# --- Tab / selection state: "increasing room size" becomes the active tab -
[void]$wsRoom.Range("I9").Select()
[void]$ws.Activate()
[void]$ws.Range("I11").Select()
